$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New value in F2
$ws.Range("F2").Value = 2

# New row 9: copy the formatting (style index "1") used by row 3's filled-in
# cells onto C9:D9, then set the actual values for the row.
$ws.Range("A3").Copy()
$ws.Range("C9:D9").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A9").Value = "ez"
$ws.Range("C9").Value = "numpy.random"
$ws.Range("D9").Value = "normal"
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = 1

# Update the active selection to match the new cursor position
$ws.Range("F14").Select()
